$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.434.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.828.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.50%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.75%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.81%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4583'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.23%  '

$ws.Range("E8").Value = '  -3.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.77%  '

$ws.Range("E10").Value = '  -1.73%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9676'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.18%  '

$ws.Range("E12").Value = '  -4.11%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.879'
$ws.Range("D13").Style = "Normal"

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.828.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.095'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.70'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06608'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001024'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.94%  '

$ws.Range("E21").Value = '  -0.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.443.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.320'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.33%  '

$ws.Range("E24").Value = '  -1.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.279'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.045.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.45%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.079'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.296'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.24%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9407'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.96%  '

$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09305'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.572'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.235'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.330'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.40%  '

$ws.Range("E37").Value = '  -2.24%  '

$ws.Range("E38").Value = '  -2.66%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.165'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.99%  '

$ws.Range("E40").Value = '  -0.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.146'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5787'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1828'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.41%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.264'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5439'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.875'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06571'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '109.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.30%  '

$ws.Range("E51").Value = '  -33.91%  '
